$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup Checklist")

# Rows whose "Configurator" dropdown value was incorrectly set to
# "Blackbox Configurator" - correct them to "RotorFlight Configurator".
$rows = @(4,5,6,11,12,14,16,18,19,20,21,22,25,28,30,31,32)
foreach ($r in $rows) {
  $ws.Range("D$r").Value = "RotorFlight Configurator"
}

# Fill in the previously empty "Program Servos" row with its tool.
$ws.Range("D15").Value = "Servo programming box"

# Reset the view back to the top of the sheet.
$ws.Activate()
$ws.Range("C1").Select()
